# Updating for changing reaction rate units in wc_lang
#
# 1. "Reactions" sheet: rename the "Flux units" column header to
#    "Flux bound units" (the flux-bound columns now have their own,
#    more specific, units label).
# 2. "dFBA objectives" sheet: insert two new columns ("Reaction rate
#    units" and "Coefficient units") right before the existing
#    "Database references" column, reflecting the new units tracked
#    for dFBA objective reaction rates / coefficients.
# 3. Leave "dFBA objectives" selected/active, matching the editor's
#    final position after making the change.

$wb = $excel.ActiveWorkbook

# --- Reactions: rename "Flux units" -> "Flux bound units" ------------
$wsReactions = $wb.Worksheets.Item("Reactions")
$wsReactions.Activate()
$wsReactions.Range("H1").Value = "Flux bound units"
$wsReactions.Range("H1").Select()

# --- dFBA objectives: insert 2 columns with new unit headers ---------
$wsObjectives = $wb.Worksheets.Item("dFBA objectives")
$wsObjectives.Activate()
$wsObjectives.Range("F1:G1").EntireColumn.Insert()
$wsObjectives.Range("F1").Value = "Reaction rate units"
$wsObjectives.Range("G1").Value = "Coefficient units"
$wsObjectives.Range("F1").Select()
